# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# header style used by the existing header row (e.g. H1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): copy the formatting of the last existing header
# cell (H1) onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-15
$data = @{
    2  = @(7, 8)
    3  = @(9, 9)
    4  = @(6, 7)
    5  = @(8, 8)
    6  = @(8, 8)
    7  = @(8, 8)
    8  = @(8, 9)
    9  = @(6, 6)
    10 = @(4, 4)
    11 = @(8, 8)
    12 = @(6, 7)
    13 = @(5, 5)
    14 = @(6, 6)
    15 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
